$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header: I1 changes from "Reason" to "Food Items"; column J removed entirely
$ws.Range("I1").Value = "Food Items"
$ws.Range("J1").Value = $null

# Delete old row 2 data in column J (Food Items column moved into I, J cleared)
$ws.Range("J2").Value = $null

# Propagate the existing date style (s="1") from B2 down to B3:B4 by copying
# formats, so the same shared style index is reused instead of new ones being
# allocated.
$ws.Range("B2").Copy()
$ws.Range("B3:B4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 2 - new data (Bill No 163)
$ws.Range("A2").Value = 163
$ws.Range("B2").Value = 45704.22928240741
$ws.Range("C2").Value = "Ajay Francis Anchan"
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 40
$ws.Range("F2").Value = 0.4
$ws.Range("G2").Value = 0.4
$ws.Range("H2").Value = 0.8
$ws.Range("I2").Value = "Belgian Coffee (x1)"

# Row 3 - new row (Bill No 162)
$ws.Range("A3").Value = 162
$ws.Range("B3").Value = 45704.22928240741
$ws.Range("C3").Value = "Ajay Francis Anchan"
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 80
$ws.Range("F3").Value = 0.8
$ws.Range("G3").Value = 0.8
$ws.Range("H3").Value = 1.6
$ws.Range("I3").Value = "8PM Coffee (x1), Coffee Italia (x1)"

# Row 4 - new row (Bill No 160)
$ws.Range("A4").Value = 160
$ws.Range("B4").Value = 45704.22928240741
$ws.Range("C4").Value = "Ajay Francis Anchan"
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 120
$ws.Range("F4").Value = 2.73
$ws.Range("G4").Value = 2.73
$ws.Range("H4").Value = 5.45
$ws.Range("I4").Value = "Veg Wrap (x1), Chicken Schezwan Wrap (x1)"

# Clear the now-unused column J entirely
$ws.Columns.Item(10).Clear()
